# Update "想去人数" (F column) counts for several events on the
# "展览" and "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 7730
$wsExpo.Range("F5").Value = 7820
$wsExpo.Range("F9").Value = 447
$wsExpo.Range("F15").Value = 306
$wsExpo.Range("F19").Value = 394
$wsExpo.Range("F20").Value = 146
$wsExpo.Range("F23").Value = 610
$wsExpo.Range("F24").Value = 2199
$wsExpo.Range("F25").Value = 731
$wsExpo.Range("F29").Value = 610

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 7730
$wsAll.Range("F7").Value = 7820
$wsAll.Range("F11").Value = 447
$wsAll.Range("F21").Value = 306
$wsAll.Range("F28").Value = 394
$wsAll.Range("F29").Value = 146
$wsAll.Range("F32").Value = 610
$wsAll.Range("F33").Value = 2199
$wsAll.Range("F34").Value = 731
$wsAll.Range("F39").Value = 610
